$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 749.125
$ws.Cells.Item(28, 10).Value = 1499.2
$ws.Cells.Item(28, 12).Value = 1499.2
$ws.Cells.Item(28, 14).Value = -2469.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(69, 8).Value = 7710.7144
$ws.Cells.Item(69, 9).Value = 7000
$ws.Cells.Item(69, 10).Value = 8658.333000000001
$ws.Cells.Item(69, 11).Value = 21000
$ws.Cells.Item(69, 12).Value = 25974.999
$ws.Cells.Item(69, 13).Value = -20126
$ws.Cells.Item(69, 14).Value = -27722.999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(72, 8).Value = 7710.7144
$ws.Cells.Item(72, 9).Value = 7000
$ws.Cells.Item(72, 10).Value = 8658.333000000001
$ws.Cells.Item(72, 11).Value = 63000
$ws.Cells.Item(72, 12).Value = 77924.997
$ws.Cells.Item(72, 13).Value = -58632
$ws.Cells.Item(72, 14).Value = -86660.997

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 5989.5884
$ws.Cells.Item(86, 9).Value = 9348.799999999999
$ws.Cells.Item(86, 10).Value = 4589.9165
$ws.Cells.Item(86, 11).Value = 9348.799999999999
$ws.Cells.Item(86, 12).Value = 4589.9165
$ws.Cells.Item(86, 13).Value = -8225.799999999999
$ws.Cells.Item(86, 14).Value = -6835.9165

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 5989.5884
$ws.Cells.Item(89, 9).Value = 9348.799999999999
$ws.Cells.Item(89, 10).Value = 4589.9165
$ws.Cells.Item(89, 11).Value = 46744
$ws.Cells.Item(89, 12).Value = 22949.5825
$ws.Cells.Item(89, 13).Value = -41128
$ws.Cells.Item(89, 14).Value = -34181.5825

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2011.6666
$ws.Cells.Item(100, 10).Value = 2372.8572
$ws.Cells.Item(100, 12).Value = 2372.8572
$ws.Cells.Item(100, 14).Value = -3454.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1354.42
$ws.Cells.Item(137, 9).Value = 844.41174
$ws.Cells.Item(137, 10).Value = 1617.1515
$ws.Cells.Item(137, 11).Value = 2533.23522
$ws.Cells.Item(137, 12).Value = 4851.4545
$ws.Cells.Item(137, 13).Value = 16.76477999999997
$ws.Cells.Item(137, 14).Value = -9951.4545

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2683.3108
$ws.Cells.Item(138, 9).Value = 1382
$ws.Cells.Item(138, 10).Value = 3347.2449
$ws.Cells.Item(138, 11).Value = 4146
$ws.Cells.Item(138, 12).Value = 10041.7347
$ws.Cells.Item(138, 13).Value = 994
$ws.Cells.Item(138, 14).Value = -20321.7347

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(37, 8).Value = 13462.667
$ws.Cells.Item(37, 10).Value = 20022.8
$ws.Cells.Item(37, 12).Value = 20022.8
$ws.Cells.Item(37, 14).Value = -20568.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 12946.667
$ws.Cells.Item(44, 10).Value = 12946.667
$ws.Cells.Item(44, 12).Value = 12946.667
$ws.Cells.Item(44, 14).Value = -13922.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(55, 8).Value = 14320
$ws.Cells.Item(55, 10).Value = 14985
$ws.Cells.Item(55, 12).Value = 14985
$ws.Cells.Item(55, 14).Value = -15615

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 3300
$ws.Cells.Item(63, 9).Value = 2933.3333
$ws.Cells.Item(63, 10).Value = 3437.5
$ws.Cells.Item(63, 11).Value = 2933.3333
$ws.Cells.Item(63, 12).Value = 3437.5
$ws.Cells.Item(63, 13).Value = -2247.3333
$ws.Cells.Item(63, 14).Value = -4809.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 3300
$ws.Cells.Item(66, 9).Value = 2933.3333
$ws.Cells.Item(66, 10).Value = 3437.5
$ws.Cells.Item(66, 11).Value = 14666.6665
$ws.Cells.Item(66, 12).Value = 17187.5
$ws.Cells.Item(66, 13).Value = -11234.6665
$ws.Cells.Item(66, 14).Value = -24051.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 24729.285
$ws.Cells.Item(80, 9).Value = 9999
$ws.Cells.Item(80, 10).Value = 27184.334
$ws.Cells.Item(80, 11).Value = 9999
$ws.Cells.Item(80, 12).Value = 27184.334
$ws.Cells.Item(80, 13).Value = -9001
$ws.Cells.Item(80, 14).Value = -29180.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(83, 8).Value = 24729.285
$ws.Cells.Item(83, 9).Value = 9999
$ws.Cells.Item(83, 10).Value = 27184.334
$ws.Cells.Item(83, 11).Value = 29997
$ws.Cells.Item(83, 12).Value = 81553.00199999999
$ws.Cells.Item(83, 13).Value = -25005
$ws.Cells.Item(83, 14).Value = -91537.00199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 18613.428
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 18613.428
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 18613.428
$ws.Cells.Item(35, 13).ClearContents()
$ws.Cells.Item(35, 14).Value = -19233.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 23235.666
$ws.Cells.Item(82, 9).Value = 2833
$ws.Cells.Item(82, 10).Value = 33437
$ws.Cells.Item(82, 11).Value = 2833
$ws.Cells.Item(82, 12).Value = 33437
$ws.Cells.Item(82, 13).Value = -2450
$ws.Cells.Item(82, 14).Value = -34203

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(85, 8).Value = 23235.666
$ws.Cells.Item(85, 9).Value = 2833
$ws.Cells.Item(85, 10).Value = 33437
$ws.Cells.Item(85, 11).Value = 2833
$ws.Cells.Item(85, 12).Value = 33437
$ws.Cells.Item(85, 13).Value = -1507
$ws.Cells.Item(85, 14).Value = -36089

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 78782.16
$ws.Cells.Item(105, 9).Value = 92661.55
$ws.Cells.Item(105, 10).Value = 68603.92999999999
$ws.Cells.Item(105, 11).Value = 92661.55
$ws.Cells.Item(105, 12).Value = 68603.92999999999
$ws.Cells.Item(105, 13).Value = -90914.55
$ws.Cells.Item(105, 14).Value = -72097.92999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 20835664
$ws.Cells.Item(107, 9).Value = 37038200
$ws.Cells.Item(107, 10).Value = 3828.4285
$ws.Cells.Item(107, 11).Value = 37038200
$ws.Cells.Item(107, 12).Value = 3828.4285
$ws.Cells.Item(107, 13).Value = -37036280
$ws.Cells.Item(107, 14).Value = -7668.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1661.1052
$ws.Cells.Item(132, 9).Value = 1204.4546
$ws.Cells.Item(132, 10).Value = 2289
$ws.Cells.Item(132, 11).Value = 3613.3638
$ws.Cells.Item(132, 12).Value = 6867
$ws.Cells.Item(132, 13).Value = -1083.3638
$ws.Cells.Item(132, 14).Value = -11927

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 5539143
$ws.Cells.Item(4, 9).Value = 8307714
$ws.Cells.Item(4, 11).Value = 24923142
$ws.Cells.Item(4, 13).Value = -24923030

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 1663.4
$ws.Cells.Item(17, 9).Value = 1533.4166
$ws.Cells.Item(17, 10).Value = 1858.375
$ws.Cells.Item(17, 11).Value = 4600.2498
$ws.Cells.Item(17, 12).Value = 5575.125
$ws.Cells.Item(17, 13).Value = -4431.2498
$ws.Cells.Item(17, 14).Value = -5913.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 2267.9412
$ws.Cells.Item(139, 9).Value = 922.6923
$ws.Cells.Item(139, 11).Value = 2768.0769
$ws.Cells.Item(139, 13).Value = 2371.9231

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3996
$ws.Cells.Item(80, 9).Value = 8800
$ws.Cells.Item(80, 11).Value = 8800
$ws.Cells.Item(80, 13).Value = -7802

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3996
$ws.Cells.Item(83, 9).Value = 8800
$ws.Cells.Item(83, 11).Value = 44000
$ws.Cells.Item(83, 13).Value = -39008

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3486.3076
$ws.Cells.Item(132, 9).Value = 3032.4
$ws.Cells.Item(132, 10).Value = 4999.3335
$ws.Cells.Item(132, 11).Value = 9097.200000000001
$ws.Cells.Item(132, 12).Value = 14998.0005
$ws.Cells.Item(132, 13).Value = -6567.200000000001
$ws.Cells.Item(132, 14).Value = -20058.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(134, 8).Value = 12371.286
$ws.Cells.Item(134, 10).Value = 12371.286
$ws.Cells.Item(134, 12).Value = 37113.858
$ws.Cells.Item(134, 14).Value = -42183.858

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4771
$ws.Cells.Item(132, 9).Value = 5706.625
$ws.Cells.Item(132, 10).Value = 3619.4614
$ws.Cells.Item(132, 11).Value = 17119.875
$ws.Cells.Item(132, 12).Value = 10858.3842
$ws.Cells.Item(132, 13).Value = -14589.875
$ws.Cells.Item(132, 14).Value = -15918.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 743.6667
$ws.Cells.Item(113, 9).Value = 383.33334
$ws.Cells.Item(113, 11).Value = 1150.00002
$ws.Cells.Item(113, 13).Value = 1019.99998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1729.2122
$ws.Cells.Item(132, 9).Value = 1323.5264
$ws.Cells.Item(132, 10).Value = 2279.7856
$ws.Cells.Item(132, 11).Value = 3970.5792
$ws.Cells.Item(132, 12).Value = 6839.3568
$ws.Cells.Item(132, 13).Value = -1440.5792
$ws.Cells.Item(132, 14).Value = -11899.3568

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1382.5
$ws.Cells.Item(136, 9).Value = 698.8333
$ws.Cells.Item(136, 10).Value = 2749.8333
$ws.Cells.Item(136, 11).Value = 2096.4999
$ws.Cells.Item(136, 12).Value = 8249.499899999999
$ws.Cells.Item(136, 13).Value = 453.5001000000002
$ws.Cells.Item(136, 14).Value = -13349.4999
